$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newBody = "Hello {0},`n<br><br>`nwe will get back to soon as per your mail`n<br><br>`nThank you! <br>`nChethan P<br>"

$ws.Range("B2").Value = $newBody
